$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Append an "order by ... limit 100" clause to each of the three Cypher
# queries stored in B2 (CasesTab), B3 (SamplesTab) and B4 (FilesTab).

$casesQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"

$samplesQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"

$filesQuery = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = $filesQuery + "`n order By f.file_name ASC LIMIT 100"

# Update the active selection to match the last edited cell.
$ws.Range("B4").Select() | Out-Null
